# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 287
$ws1.Range("F3").Value = 1173
$ws1.Range("F4").Value = 16675
$ws1.Range("F5").Value = 25
$ws1.Range("F6").Value = 1633
$ws1.Range("F8").Value = 2
$ws1.Range("F9").Value = 363
$ws1.Range("F12").Value = 11584
$ws1.Range("F14").Value = 1267
$ws1.Range("F15").Value = 4585
$ws1.Range("F16").Value = 419
$ws1.Range("F17").Value = 402
$ws1.Range("F19").Value = 878

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 287
$ws4.Range("F4").Value = 1173
$ws4.Range("F5").Value = 16675
$ws4.Range("F6").Value = 25
$ws4.Range("F7").Value = 1633
$ws4.Range("F9").Value = 2
$ws4.Range("F10").Value = 363
$ws4.Range("F15").Value = 11584
$ws4.Range("F17").Value = 1267
$ws4.Range("F18").Value = 4585
$ws4.Range("F19").Value = 419
$ws4.Range("F20").Value = 402
$ws4.Range("F22").Value = 878
